$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new VOC id, new correlation values, new influence label
$ws.Range("A2").Value = "000100-41-4_group4"
$ws.Range("B2").Value = -0.082
$ws.Range("C2").Value = -0.076
$ws.Range("D2").Value = 0.037
$ws.Range("E2").Value = "Unlikely"

# Row 3: VOC id (A3) and label (E3) stay the same, only the three correlation values change
$ws.Range("B3").Value = 0.08599999999999999
$ws.Range("C3").Value = -0.004
$ws.Range("D3").Value = 0.004

# Row 4: new VOC id, new correlation values; label (E4) stays "Unlikely"
$ws.Range("A4").Value = "000124-19-6_group1"
$ws.Range("B4").Value = 0.07099999999999999
$ws.Range("C4").Value = 0.008
$ws.Range("D4").Value = 0.01

# Row 5: new VOC id, new correlation values; label (E5) stays "Unlikely"
$ws.Range("A5").Value = "000098-86-2_group2"
$ws.Range("B5").Value = 0.098
$ws.Range("C5").Value = 0.061
$ws.Range("D5").Value = 0.002

# Row 6 is removed entirely, shrinking the used range to A1:E5
$ws.Rows.Item(6).Delete()
